$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213302731513977
$ws.Range("B1").Value = 2.584191083908081
$ws.Range("C1").Value = 4.342037677764893
$ws.Range("D1").Value = 2.073843002319336
$ws.Range("E1").Value = 1.167407155036926
